$wb = $excel.ActiveWorkbook

# ALC!row18
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 569
$ws.Cells.Item(18, 9).Value = 564.2857
$ws.Cells.Item(18, 10).Value = 602
$ws.Cells.Item(18, 11).Value = 564.2857
$ws.Cells.Item(18, 12).Value = 602
$ws.Cells.Item(18, 13).Value = -280.2857
$ws.Cells.Item(18, 14).Value = -1170

# ALC!row40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 4344.433
$ws.Cells.Item(40, 9).Value = 2276.1538
$ws.Cells.Item(40, 10).Value = 5926.0586
$ws.Cells.Item(40, 11).Value = 2276.1538
$ws.Cells.Item(40, 12).Value = 5926.0586
$ws.Cells.Item(40, 13).Value = -2101.1538
$ws.Cells.Item(40, 14).Value = -6276.0586

# ALC!row55
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 154.58333
$ws.Cells.Item(55, 9).Value = 141.36363
$ws.Cells.Item(55, 11).Value = 141.36363
$ws.Cells.Item(55, 13).Value = 72.63637

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 4795.5625
$ws.Cells.Item(132, 9).Value = 4781.933
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 14345.799
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -11815.799
$ws.Cells.Item(132, 14).Value = -20060

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5024.4805
$ws.Cells.Item(32, 9).Value = 1866.6
$ws.Cells.Item(32, 11).Value = 1866.6
$ws.Cells.Item(32, 13).Value = -1579.6

# ARM!row63
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 2767.647
$ws.Cells.Item(63, 9).Value = 1913.6364
$ws.Cells.Item(63, 10).Value = 4333.3335
$ws.Cells.Item(63, 11).Value = 1913.6364
$ws.Cells.Item(63, 12).Value = 4333.3335
$ws.Cells.Item(63, 13).Value = -1227.6364
$ws.Cells.Item(63, 14).Value = -5705.3335

# ARM!row66
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 2767.647
$ws.Cells.Item(66, 9).Value = 1913.6364
$ws.Cells.Item(66, 10).Value = 4333.3335
$ws.Cells.Item(66, 11).Value = 9568.182000000001
$ws.Cells.Item(66, 12).Value = 21666.6675
$ws.Cells.Item(66, 13).Value = -6136.182000000001
$ws.Cells.Item(66, 14).Value = -28530.6675

# ARM!row88
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 2451.125
$ws.Cells.Item(88, 10).Value = 2451.125
$ws.Cells.Item(88, 12).Value = 2451.125
$ws.Cells.Item(88, 14).Value = -3263.125

# ARM!row91
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 2451.125
$ws.Cells.Item(91, 10).Value = 2451.125
$ws.Cells.Item(91, 12).Value = 2451.125
$ws.Cells.Item(91, 14).Value = -5259.125

# BSM!row22
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 613.7406999999999
$ws.Cells.Item(22, 9).Value = 521.35297
$ws.Cells.Item(22, 11).Value = 521.35297
$ws.Cells.Item(22, 13).Value = -348.35297

# BSM!row80
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 58049.57
$ws.Cells.Item(80, 10).Value = 21269.2
$ws.Cells.Item(80, 12).Value = 21269.2
$ws.Cells.Item(80, 14).Value = -23265.2

# BSM!row83
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value = 58049.57
$ws.Cells.Item(83, 10).Value = 21269.2
$ws.Cells.Item(83, 12).Value = 106346
$ws.Cells.Item(83, 14).Value = -116330

# CRP!row50
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 40950.168
$ws.Cells.Item(50, 10).Value = 56424.5
$ws.Cells.Item(50, 12).Value = 56424.5
$ws.Cells.Item(50, 14).Value = -57674.5

# CRP!row51
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(51, 8).Value = 53657.5
$ws.Cells.Item(51, 10).Value = 53657.5
$ws.Cells.Item(51, 12).Value = 53657.5
$ws.Cells.Item(51, 14).Value = -55129.5

# CRP!row60
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(60, 8).Value = 21119.5
$ws.Cells.Item(60, 9).Value = 11719.6
$ws.Cells.Item(60, 10).Value = 26341.666
$ws.Cells.Item(60, 11).Value = 11719.6
$ws.Cells.Item(60, 12).Value = 26341.666
$ws.Cells.Item(60, 13).Value = -11208.6
$ws.Cells.Item(60, 14).Value = -27363.666

# CRP!row61
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(61, 8).Value = 53657.5
$ws.Cells.Item(61, 10).Value = 53657.5
$ws.Cells.Item(61, 12).Value = 53657.5
$ws.Cells.Item(61, 14).Value = -54353.5

# CUL!row3
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 2630.7
$ws.Cells.Item(3, 9).Value = 2106
$ws.Cells.Item(3, 11).Value = 6318
$ws.Cells.Item(3, 13).Value = -6206

# CUL!row59
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(59, 8).Value = 7994
$ws.Cells.Item(59, 9).Value = 7994
$ws.Cells.Item(59, 11).Value = 23982
$ws.Cells.Item(59, 13).Value = -23442

# CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 967
$ws.Cells.Item(122, 10).Value = 950.5
$ws.Cells.Item(122, 12).Value = 8554.5
$ws.Cells.Item(122, 14).Value = -13454.5

# CUL!row133
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(133, 8).Value = 8972.5
$ws.Cells.Item(133, 9).Value = 3677.3635
$ws.Cells.Item(133, 10).Value = 15444.333
$ws.Cells.Item(133, 11).Value = 11032.0905
$ws.Cells.Item(133, 12).Value = 46332.999
$ws.Cells.Item(133, 13).Value = -5972.0905
$ws.Cells.Item(133, 14).Value = -56452.999

# CUL!row137
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 6108.6
$ws.Cells.Item(137, 9).Value = 2591.3333
$ws.Cells.Item(137, 10).Value = 7616
$ws.Cells.Item(137, 11).Value = 7773.999899999999
$ws.Cells.Item(137, 12).Value = 22848
$ws.Cells.Item(137, 13).Value = -2673.999899999999
$ws.Cells.Item(137, 14).Value = -33048

# CUL!row139
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 2064.875
$ws.Cells.Item(139, 9).Value = 1484
$ws.Cells.Item(139, 11).Value = 4452
$ws.Cells.Item(139, 13).Value = 688

# GSM!row5
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).ClearContents()

# GSM!row9
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 3295.875
$ws.Cells.Item(9, 9).Value = 3648.3333
$ws.Cells.Item(9, 11).Value = 3648.3333
$ws.Cells.Item(9, 13).Value = -3478.3333

# GSM!row80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 10675.069
$ws.Cells.Item(80, 9).Value = 9482.846
$ws.Cells.Item(80, 11).Value = 9482.846
$ws.Cells.Item(80, 13).Value = -8484.846

# GSM!row83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 10675.069
$ws.Cells.Item(83, 9).Value = 9482.846
$ws.Cells.Item(83, 11).Value = 47414.23
$ws.Cells.Item(83, 13).Value = -42422.23

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 7109
$ws.Cells.Item(122, 9).Value = 5205.067
$ws.Cells.Item(122, 10).Value = 16628.666
$ws.Cells.Item(122, 11).Value = 15615.201
$ws.Cells.Item(122, 12).Value = 49885.99800000001
$ws.Cells.Item(122, 13).Value = -13165.201
$ws.Cells.Item(122, 14).Value = -54785.99800000001

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 7034.544
$ws.Cells.Item(132, 9).Value = 6871.82
$ws.Cells.Item(132, 11).Value = 20615.46
$ws.Cells.Item(132, 13).Value = -18085.46

# LTW!row9
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 8052.857
$ws.Cells.Item(9, 9).Value = 9061.666999999999
$ws.Cells.Item(9, 11).Value = 9061.666999999999
$ws.Cells.Item(9, 13).Value = -8837.666999999999

# LTW!row30
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(30, 8).Value = 4000
$ws.Cells.Item(30, 9).Value = 4000
$ws.Cells.Item(30, 11).Value = 4000
$ws.Cells.Item(30, 13).Value = -3892

# LTW!row35
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(35, 8).Value = 7166.3335
$ws.Cells.Item(35, 10).Value = 9499
$ws.Cells.Item(35, 12).Value = 9499
$ws.Cells.Item(35, 14).Value = -10171

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 5479.467
$ws.Cells.Item(136, 9).Value = 4259
$ws.Cells.Item(136, 10).Value = 7920.4
$ws.Cells.Item(136, 11).Value = 12777
$ws.Cells.Item(136, 12).Value = 23761.2
$ws.Cells.Item(136, 13).Value = -10227
$ws.Cells.Item(136, 14).Value = -28861.2

# WVR!row62
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 29990
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 29990
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 29990
$ws.Cells.Item(62, 14).Value = -31238
$ws.Cells.Item(62, 13).ClearContents()

# WVR!row65
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 29990
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 29990
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 149950
$ws.Cells.Item(65, 14).Value = -156190
$ws.Cells.Item(65, 13).ClearContents()
